# NEVADA_2021.xlsx cleanup:
#  1. Rename header row (A1:D1) to short machine-friendly column codes.
#  2. Title-case the Spanish state/municipality names in columns A and B
#     (lowercase connector words like "de", "del", "la", "los", "el", "y"
#     get capitalized, matching a simple split-on-space .title() pass).
#  3. Drop the trailing metadata/footnote rows (1377-1381), shrinking the
#     used range back down to A1:D1375.

function ConvertTo-TitleCaseEs($s) {
    if ($null -eq $s -or $s.Length -eq 0) {
        return $s
    }
    $parts = $s.Split(" ")
    $out = @()
    foreach ($p in $parts) {
        if ($p.Length -gt 0) {
            # NB: use [string]::Concat (not "+") - this interpreter coerces
            # "+" between two digit-looking strings (e.g. "2" and "6") into
            # *numeric* addition ("8"), which corrupts tokens like "26".
            $head = $p.Substring(0,1).ToUpper()
            $tail = $p.Substring(1)
            $out += [string]::Concat($head, $tail)
        } else {
            $out += $p
        }
    }
    return [string]::Join(" ", $out)
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header row rename -------------------------------------------------
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- 2. Title-case columns A and B for every data row ---------------------
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aText = $aCell.Text
    if ($aText.Length -gt 0) {
        $aCell.Value = ConvertTo-TitleCaseEs $aText
    }

    $bCell = $ws.Cells.Item($r, 2)
    $bText = $bCell.Text
    if ($bText.Length -gt 0) {
        $bCell.Value = ConvertTo-TitleCaseEs $bText
    }
}

# --- 3. Remove the trailing metadata / footnote rows ----------------------
$ws.Range("A1377:A1381").EntireRow.Delete()
